$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.501.45"
$ws.Range("E2").Value = "  -5.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.520.44"
$ws.Range("E3").Value = "  -5.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.67"
$ws.Range("E5").Value = "  -7.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.29"
$ws.Range("E6").Value = "  -1.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.508.09"
$ws.Range("E7").Value = "  -5.91%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.602"
$ws.Range("E8").Value = "  -5.22%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.659"
$ws.Range("E10").Value = "  -9.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.65"
$ws.Range("E11").Value = "  -7.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.142"
$ws.Range("E12").Value = "  -12.96%  "
$ws.Range("E13").Value = "  -17.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.65"
$ws.Range("E14").Value = "  -10.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.094.48"
$ws.Range("E15").Value = "  -5.04%  "
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.524.29"
$ws.Range("E17").Value = "  -5.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.07"
$ws.Range("E18").Value = "  -7.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "65.265.20"
$ws.Range("E19").Value = "  -5.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.99"
$ws.Range("E20").Value = "  -7.99%  "
$ws.Range("E21").Value = "  -8.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "387.20"
$ws.Range("E22").Value = "  -7.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.21"
$ws.Range("E23").Value = "  -11.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.41"
$ws.Range("E24").Value = "  -7.07%  "
$ws.Range("E25").Value = "  -7.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.19"
$ws.Range("E26").Value = "  -4.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.99"
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.17"
$ws.Range("E28").Value = "  -7.95%  "
$ws.Range("E29").Value = "  -9.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.79"
$ws.Range("E30").Value = "  -9.13%  "
$ws.Range("E31").Value = "  -8.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.67"
$ws.Range("E32").Value = "  -10.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.89"
$ws.Range("E33").Value = "  -5.70%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "605.91"
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "64.52"
$ws.Range("E35").Value = "  -2.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.110"
$ws.Range("E36").Value = "  -8.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "40.86"
$ws.Range("E37").Value = "  -7.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("E40").Value = "  -10.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0729"
$ws.Range("E41").Value = "  -19.00%  "
$ws.Range("E42").Value = "  -7.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.848.24"
$ws.Range("E43").Value = "  +1.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.73"
$ws.Range("E44").Value = "  -11.48%  "
$ws.Range("E45").Value = "  -9.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.37"
$ws.Range("E46").Value = "  -11.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.128"
$ws.Range("E47").Value = "  -6.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "137.26"
$ws.Range("E48").Value = "  -2.74%  "
$ws.Range("E49").Value = "  -5.68%  "
$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.44"
$ws.Range("E50").Value = "  -10.85%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.11"
$ws.Range("E51").Value = "  -12.87%  "
